$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.573375
$ws.Range("H2").Value = 28.720125
$ws.Range("I2").Value = 0.1037691388643484
$ws.Range("J2").Value = 0.1037691388643484
$ws.Range("M2").Value = 4.277890333333334
$ws.Range("N2").Value = 12.833671
$ws.Range("O2").Value = 0.04123357425337639
$ws.Range("P2").Value = 0.04123357425337638
$ws.Range("Q2").Value = 40.95384836987501
$ws.Range("R2").Value = 368.584635328875
$ws.Range("S2").Value = 0.004278772492572035
$ws.Range("T2").Value = 0.004278772492572035
$ws.Range("G3").Value = 9.573375
$ws.Range("H3").Value = 28.720125
$ws.Range("I3").Value = 0.1037691388643484
$ws.Range("J3").Value = 0.1037691388643484
$ws.Range("O3").Value = 0.4451428460610328
$ws.Range("P3").Value = 0.4451428460610327
$ws.Range("Q3").Value = 442.123025971375
$ws.Range("R3").Value = 3979.107233742375
$ws.Range("S3").Value = 0.04619208980737857
$ws.Range("T3").Value = 0.04619208980737857
$ws.Range("G4").Value = 9.573375
$ws.Range("H4").Value = 28.720125
$ws.Range("I4").Value = 0.1037691388643484
$ws.Range("J4").Value = 0.1037691388643484
$ws.Range("M4").Value = 8.558147333333332
$ws.Range("N4").Value = 25.674442
$ws.Range("O4").Value = 0.08248996024761777
$ws.Range("P4").Value = 0.08248996024761777
$ws.Range("Q4").Value = 81.93035372724999
$ws.Range("R4").Value = 737.37318354525
$ws.Range("S4").Value = 0.008559912139849628
$ws.Range("T4").Value = 0.008559912139849628
$ws.Range("G5").Value = 9.573375
$ws.Range("H5").Value = 28.720125
$ws.Range("I5").Value = 0.1037691388643484
$ws.Range("J5").Value = 0.1037691388643484
$ws.Range("M5").Value = 44.72914066666667
$ws.Range("N5").Value = 134.187422
$ws.Range("O5").Value = 0.4311336194379731
$ws.Range("P5").Value = 0.431133619437973
$ws.Range("Q5").Value = 428.20883702975
$ws.Range("R5").Value = 3853.87953326775
$ws.Range("S5").Value = 0.04473836442454816
$ws.Range("T5").Value = 0.04473836442454816
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("I6").Value = 0.5196887643218222
$ws.Range("J6").Value = 0.5196887643218222
$ws.Range("M6").Value = 4.277890333333334
$ws.Range("N6").Value = 12.833671
$ws.Range("O6").Value = 0.04123357425337639
$ws.Range("P6").Value = 0.04123357425337638
$ws.Range("Q6").Value = 205.101970455648
$ws.Range("R6").Value = 1845.917734100831
$ws.Range("S6").Value = 0.02142862525230928
$ws.Range("T6").Value = 0.02142862525230928
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("I7").Value = 0.5196887643218222
$ws.Range("J7").Value = 0.5196887643218222
$ws.Range("O7").Value = 0.4451428460610328
$ws.Range("P7").Value = 0.4451428460610327
$ws.Range("Q7").Value = 2214.207148289527
$ws.Range("S7").Value = 0.2313357356161573
$ws.Range("T7").Value = 0.2313357356161572
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("I8").Value = 0.5196887643218222
$ws.Range("J8").Value = 0.5196887643218222
$ws.Range("M8").Value = 8.558147333333332
$ws.Range("N8").Value = 25.674442
$ws.Range("O8").Value = 0.08248996024761777
$ws.Range("P8").Value = 0.08248996024761777
$ws.Range("Q8").Value = 410.3174099249736
$ws.Range("R8").Value = 3692.856689324762
$ws.Range("S8").Value = 0.04286910551004072
$ws.Range("T8").Value = 0.04286910551004072
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("I9").Value = 0.5196887643218222
$ws.Range("J9").Value = 0.5196887643218222
$ws.Range("M9").Value = 44.72914066666667
$ws.Range("N9").Value = 134.187422
$ws.Range("O9").Value = 0.4311336194379731
$ws.Range("P9").Value = 0.431133619437973
$ws.Range("Q9").Value = 2144.52315807095
$ws.Range("R9").Value = 19300.70842263855
$ws.Range("S9").Value = 0.224055297943315
$ws.Range("T9").Value = 0.2240552979433149
$ws.Range("G10").Value = 11.32006633333333
$ws.Range("H10").Value = 33.960199
$ws.Range("I10").Value = 0.122702133291269
$ws.Range("J10").Value = 0.122702133291269
$ws.Range("M10").Value = 4.277890333333334
$ws.Range("N10").Value = 12.833671
$ws.Range("O10").Value = 0.04123357425337639
$ws.Range("P10").Value = 0.04123357425337638
$ws.Range("Q10").Value = 48.42600234005879
$ws.Range("R10").Value = 435.8340210605291
$ws.Range("S10").Value = 0.005059447524113226
$ws.Range("T10").Value = 0.005059447524113225
$ws.Range("G11").Value = 11.32006633333333
$ws.Range("H11").Value = 33.960199
$ws.Range("I11").Value = 0.122702133291269
$ws.Range("J11").Value = 0.122702133291269
$ws.Range("O11").Value = 0.4451428460610328
$ws.Range("P11").Value = 0.4451428460610327
$ws.Range("Q11").Value = 522.7897143369001
$ws.Range("R11").Value = 4705.107429032101
$ws.Range("S11").Value = 0.05461997683103566
$ws.Range("T11").Value = 0.05461997683103566
$ws.Range("G12").Value = 11.32006633333333
$ws.Range("H12").Value = 33.960199
$ws.Range("I12").Value = 0.122702133291269
$ws.Range("J12").Value = 0.122702133291269
$ws.Range("M12").Value = 8.558147333333332
$ws.Range("N12").Value = 25.674442
$ws.Range("O12").Value = 0.08248996024761777
$ws.Range("P12").Value = 0.08248996024761777
$ws.Range("Q12").Value = 96.87879550377312
$ws.Range("R12").Value = 871.909159533958
$ws.Range("S12").Value = 0.01012169409749467
$ws.Range("T12").Value = 0.01012169409749467
$ws.Range("G13").Value = 11.32006633333333
$ws.Range("H13").Value = 33.960199
$ws.Range("I13").Value = 0.122702133291269
$ws.Range("J13").Value = 0.122702133291269
$ws.Range("M13").Value = 44.72914066666667
$ws.Range("N13").Value = 134.187422
$ws.Range("O13").Value = 0.4311336194379731
$ws.Range("P13").Value = 0.431133619437973
$ws.Range("Q13").Value = 506.3368393796643
$ws.Range("R13").Value = 4557.031554416978
$ws.Range("S13").Value = 0.05290101483862539
$ws.Range("T13").Value = 0.05290101483862539
$ws.Range("G14").Value = 23.41838033333333
$ws.Range("H14").Value = 70.25514099999999
$ws.Range("I14").Value = 0.2538399635225604
$ws.Range("J14").Value = 0.2538399635225604
$ws.Range("M14").Value = 4.277890333333334
$ws.Range("N14").Value = 12.833671
$ws.Range("O14").Value = 0.04123357425337639
$ws.Range("P14").Value = 0.04123357425337638
$ws.Range("Q14").Value = 100.1812628502901
$ws.Range("R14").Value = 901.6313656526111
$ws.Range("S14").Value = 0.01046672898438185
$ws.Range("T14").Value = 0.01046672898438185
$ws.Range("G15").Value = 23.41838033333333
$ws.Range("H15").Value = 70.25514099999999
$ws.Range("I15").Value = 0.2538399635225604
$ws.Range("J15").Value = 0.2538399635225604
$ws.Range("O15").Value = 0.4451428460610328
$ws.Range("P15").Value = 0.4451428460610327
$ws.Range("Q15").Value = 1081.520903163395
$ws.Range("R15").Value = 9733.688128470558
$ws.Range("S15").Value = 0.1129950438064613
$ws.Range("T15").Value = 0.1129950438064613
$ws.Range("G16").Value = 23.41838033333333
$ws.Range("H16").Value = 70.25514099999999
$ws.Range("I16").Value = 0.2538399635225604
$ws.Range("J16").Value = 0.2538399635225604
$ws.Range("M16").Value = 8.558147333333332
$ws.Range("N16").Value = 25.674442
$ws.Range("O16").Value = 0.08248996024761777
$ws.Range("P16").Value = 0.08248996024761777
$ws.Range("Q16").Value = 200.4179492007024
$ws.Range("R16").Value = 1803.761542806322
$ws.Range("S16").Value = 0.02093924850023276
$ws.Range("T16").Value = 0.02093924850023276
$ws.Range("G17").Value = 23.41838033333333
$ws.Range("H17").Value = 70.25514099999999
$ws.Range("I17").Value = 0.2538399635225604
$ws.Range("J17").Value = 0.2538399635225604
$ws.Range("M17").Value = 44.72914066666667
$ws.Range("N17").Value = 134.187422
$ws.Range("O17").Value = 0.4311336194379731
$ws.Range("P17").Value = 0.431133619437973
$ws.Range("Q17").Value = 1047.484028115167
$ws.Range("R17").Value = 9427.356253036502
$ws.Range("S17").Value = 0.1094389422314845
$ws.Range("T17").Value = 0.1094389422314845
